$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("D1").Value = "Total"
$ws.Range("E1").Value = "Percent"

# Copy header style (bold, centered, bordered) from C1 to D1:E1
$ws.Range("C1").Copy()
$ws.Range("D1:E1").PasteSpecial(-4122)  # xlPasteFormats

# Force percent column to be text so literal "%" values are kept as strings
$ws.Range("E2:E5").NumberFormat = "@"

# Data rows
$ws.Range("D2").Value = 1419
$ws.Range("E2").Value = "18.94%"

$ws.Range("D3").Value = 902
$ws.Range("E3").Value = "12.04%"

$ws.Range("D4").Value = 2057
$ws.Range("E4").Value = "27.46%"

$ws.Range("D5").Value = 3114
$ws.Range("E5").Value = "41.56%"

# Copy data style (centered, General number format) from C2:C5 to D2:E5
# (applied after setting values so the stored text type is preserved while
# the number format / style index matches column C's style)
$ws.Range("C2:C5").Copy()
$ws.Range("D2:E5").PasteSpecial(-4122)  # xlPasteFormats

$excel.CutCopyMode = 0
